# Rename the worksheet from "Property1" to "DataNode" to unify the
# conception of DataNode / DataTable / Entity across the data config sheets.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Move the active selection, matching the cursor position left behind by
# the author when they saved the workbook.
$ws.Range("C41").Select()
